$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Rows 235-239: re-sync odds data rows that were re-ordered in source feed ----
# Row 235
$ws.Range("B235").Value = 6861095
$ws.Range("F235").Value = "FC Botosani"
$ws.Range("G235").Value = "Farul Constanta"
$ws.Range("H235").Value = 0
$ws.Range("J235").Value = "D"
$ws.Range("K235").Value = 3.75
$ws.Range("L235").Value = 3.4
$ws.Range("M235").Value = 1.909
$ws.Range("N235").Value = 3.1
$ws.Range("P235").Value = 2.375
$ws.Range("Q235").Value = 0.25
$ws.Range("R235").Value = 1.775
$ws.Range("S235").Value = 2.1
$ws.Range("T235").Value = 2
$ws.Range("U235").Value = 1.8
$ws.Range("V235").Value = 2.05
$ws.Range("W235").Value = -1
$ws.Range("X235").Value = 2
$ws.Range("Z235").Value = 0.3875
$ws.Range("AA235").Value = -0.5
$ws.Range("AC235").Value = 1.05

# Row 236
$ws.Range("B236").Value = 6865915
$ws.Range("F236").Value = "FC Voluntari"
$ws.Range("G236").Value = "Universitatea Cluj"
$ws.Range("H236").Value = 0
$ws.Range("J236").Value = "D"
$ws.Range("K236").Value = 3.5
$ws.Range("L236").Value = 3.25
$ws.Range("M236").Value = 2.05
$ws.Range("N236").Value = 3.4
$ws.Range("O236").Value = 3.1
$ws.Range("P236").Value = 2.15
$ws.Range("Q236").Value = 0.25
$ws.Range("R236").Value = 1.975
$ws.Range("S236").Value = 1.875
$ws.Range("U236").Value = 2.05
$ws.Range("V236").Value = 1.75
$ws.Range("W236").Value = -1
$ws.Range("X236").Value = 2.1
$ws.Range("Z236").Value = 0.4875
$ws.Range("AC236").Value = 0.75

# Row 237
$ws.Range("B237").Value = 6836277
$ws.Range("F237").Value = "CFR Cluj"
$ws.Range("G237").Value = "AFC Hermannstadt"
$ws.Range("I237").Value = 0
$ws.Range("J237").Value = "H"
$ws.Range("K237").Value = 1.7
$ws.Range("L237").Value = 3.4
$ws.Range("M237").Value = 5
$ws.Range("N237").Value = 1.65
$ws.Range("O237").Value = 3.5
$ws.Range("P237").Value = 5.25
$ws.Range("Q237").Value = -0.75
$ws.Range("W237").Value = 0.6499999999999999
$ws.Range("Y237").Value = -1
$ws.Range("Z237").Value = 0.425
$ws.Range("AA237").Value = -0.5
$ws.Range("AB237").Value = -1
$ws.Range("AC237").Value = 0.9750000000000001

# Row 238
$ws.Range("B238").Value = 6852370
$ws.Range("F238").Value = "Dinamo Bucharest"
$ws.Range("G238").Value = "ACS UTA Batrana Doamna"
$ws.Range("H238").Value = 1
$ws.Range("J238").Value = "H"
$ws.Range("K238").Value = 2.55
$ws.Range("L238").Value = 2.875
$ws.Range("M238").Value = 3
$ws.Range("N238").Value = 2.375
$ws.Range("O238").Value = 3
$ws.Range("P238").Value = 3.1
$ws.Range("Q238").Value = -0.25
$ws.Range("R238").Value = 2
$ws.Range("S238").Value = 1.85
$ws.Range("U238").Value = 1.975
$ws.Range("V238").Value = 1.875
$ws.Range("W238").Value = 1.375
$ws.Range("X238").Value = -1
$ws.Range("Z238").Value = 1
$ws.Range("AA238").Value = -1
$ws.Range("AC238").Value = 0.875

# Row 239
$ws.Range("B239").Value = 6870268
$ws.Range("F239").Value = "Petrolul Ploiesti"
$ws.Range("G239").Value = "ACS Sepsi"
$ws.Range("H239").Value = 1
$ws.Range("I239").Value = 2
$ws.Range("J239").Value = "A"
$ws.Range("K239").Value = 2.8
$ws.Range("L239").Value = 3
$ws.Range("M239").Value = 2.55
$ws.Range("N239").Value = 3
$ws.Range("O239").Value = 3.2
$ws.Range("P239").Value = 2.3
$ws.Range("R239").Value = 1.85
$ws.Range("S239").Value = 2
$ws.Range("T239").Value = 2.25
$ws.Range("U239").Value = 1.875
$ws.Range("V239").Value = 1.975
$ws.Range("X239").Value = -1
$ws.Range("Y239").Value = 1.3
$ws.Range("Z239").Value = -1
$ws.Range("AA239").Value = 1
$ws.Range("AB239").Value = 0.875
$ws.Range("AC239").Value = -1

# ---- Row 261: updated pre-match odds ----
$ws.Range("N261").Value = 2.15
$ws.Range("O261").Value = 3.25
$ws.Range("P261").Value = 3.5
$ws.Range("Q261").Value = -0.25
$ws.Range("R261").Value = 1.875
$ws.Range("S261").Value = 1.975
